$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 (shifts existing rows 7.. down by one),
# matching the FxE matrix "output" configuration row that was added
# ahead of the existing "output_efficiency" row.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new FxE "output" entry.
$ws.Range("A7").Value = "CHE"
$ws.Range("B7").Value = "ext_fueloil"
$ws.Range("C7").Value = "output"
$ws.Range("D7").Value = "configuration_fxe"
$ws.Range("F7").Value = "fueloil"
$ws.Range("G7").Value = 1

# The autofilter / filter-database range grew by one row because of the
# inserted row; refresh it to cover A5:L573 instead of A5:L572.
$ws.AutoFilterMode = $false
[void]$ws.Range("A5:L573").AutoFilter()

$filterName = $wb.Names.Item("Sheet1!_FilterDatabase")
$filterName.RefersTo = "=Sheet1!`$A`$5:`$L`$573"

# Restore the selection to the cell that was active after the edit.
[void]$ws.Range("E7").Select()
